# SI4825 KIT TH BOM - "Fixed issues discovered in first build" edit
#
# - C9 (AC coupling capacitor) changed value to 22uF and split out of the
#   "C3-C9" group into its own BOM line (group becomes "C3-C6, C8").
# - R8 (jumper to bypass op amp) marked as NOSTUFF (not populated).
# - U2 (op amp) changed from SMD MCP6002T-I/SN to through-hole MCP6002-E/P.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row for C9 right after the C3-C9 row (row 3), pushing
# everything below down by one.
$ws.Rows("4:4").Insert()

# Row 3: "C3-C9" group loses C9, qty drops from 7 to 5.
$ws.Range("A3").Value = "C3-C6, C8"
$ws.Range("E3").Value = 5

# Row 4 (new): C9, now a 22uF 5mm radial X5R capacitor from TDK.
$ws.Range("G4").Value = "http://www.digikey.com/product-detail/en/tdk-corporation/FK24X5R0J226M/445-8492-ND/2815422"
$ws.Range("A4").Value = "C9"
$ws.Range("B4").Value = "22uF 5mm radial Capacitor X5R"
$ws.Range("C4").Value = "TDK"
$ws.Range("D4").Value = "FK24X5R0J226M"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.319

# Row 14 (was row 13 before the insert): R8 jumper is no longer stuffed.
$ws.Range("E14").Value = "NOSTUFF"

# Row 18 (was row 17 before the insert): U2 swapped for the through-hole
# PDIP-8 op amp.
$ws.Range("G18").Value = "http://www.digikey.com/product-detail/en/microchip-technology/MCP6002-E-P/MCP6002-E-P-ND/683196"
$ws.Range("D18").Value = "MCP6002-E/P"
$ws.Range("B18").Value = "Dual low-voltage op amp (PDIP-8)"
$ws.Range("C18").Value = "Microchip"
$ws.Range("F18").Value = 0.3

# Restore the active selection to G4, matching the author's saved view.
$ws.Range("G4").Select()
